$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; B=3.272327238179451;  C=1.626987699542094;  D=0.7210945179870265; E=0.5333859586016987; G=6.15379541431027}
    @{Row=3; B=1.445647641019636;  C=1.626987699542094;  D=0.1496068669990043; E=0.5333859586016987; G=3.755628166162433}
    @{Row=4; B=1.445647641019636;  C=1.626987699542094;  D=0.1496068669990043; E=0.5333859586016987; G=3.755628166162433}
    @{Row=5; B=1.445647641019636;  C=1.626987699542094;  D=0.1496068669990043; E=0.5333859586016987; G=3.755628166162433}
    @{Row=6; B=0.1169995834814548; C=0.002658071450198252;D=0.7210945179870265; E=0.5333859586016987; G=1.374138131520378}
    @{Row=7; B=3.272327238179451;  C=1.626987699542094;  D=0.1496068669990043; E=0.5333859586016987; G=5.582307763322248}
    @{Row=8; B=0.04172184405617529;C=0.04103571897497393;D=0.7210945179870265; E=0.5333859586016987; G=1.337238039619874}
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item($r, 7).Value = $entry.G
}
